$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header (text content unchanged, rewritten for completeness)
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Prognose"
$ws.Range("C1").Value = "surveys"
$ws.Range("D1").Value = "production"
$ws.Range("E1").Value = "orders"
$ws.Range("F1").Value = "turnover"
$ws.Range("G1").Value = "financial"
$ws.Range("H1").Value = "labor market"
$ws.Range("I1").Value = "prices"
$ws.Range("J1").Value = "national accounts"
$ws.Range("K1").Value = "Revision"

# Rows 2-7 column A: new rolling date labels (must remain text, not auto-converted to Excel dates)
$dateCells = "A2","A3","A4","A5","A6","A7"
foreach ($addr in $dateCells) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("A2").Value = "2025-09-30"
$ws.Range("A3").Value = "2025-10-15"
$ws.Range("A4").Value = "2025-10-30"
$ws.Range("A5").Value = "2025-11-15"
$ws.Range("A6").Value = "2025-11-30"
$ws.Range("A7").Value = "2025-12-15"
foreach ($addr in $dateCells) { $ws.Range($addr).ClearFormats() }

# Rows 2-7 columns B:K: updated nowcast / revision figures
$ws.Range("B2").Value = 0.040871905345710557
$ws.Range("B3").Value = -0.015469416270583181
$ws.Range("D3").Value = -0.053050962818628304
$ws.Range("E3").Value = 0.0011137182538578287
$ws.Range("F3").Value = -0.0032751809935352739
$ws.Range("G3").Value = 0.0043413870823703019
$ws.Range("H3").Value = 0.00048207782823450542
$ws.Range("I3").Value = -0.0046669644764513411
$ws.Range("K3").Value = -0.0012853964921414607
$ws.Range("B4").Value = 0.28366676301632843
$ws.Range("C4").Value = 0.24291870380191261
$ws.Range("E4").Value = 0.00028655976289332239
$ws.Range("F4").Value = 0.00026703430868886178
$ws.Range("H4").Value = -0.0068812128813553076
$ws.Range("I4").Value = -0.012117467371370875
$ws.Range("J4").Value = 0.077931407356559929
$ws.Range("K4").Value = -0.0032688456904169105
$ws.Range("B5").Value = 0.38631825920457591
$ws.Range("D5").Value = 0.12217228098294355
$ws.Range("E5").Value = 0.025952390243862408
$ws.Range("F5").Value = -0.051540676641486212
$ws.Range("G5").Value = -0.0027566232280068323
$ws.Range("H5").Value = -0.0066278655185185885
$ws.Range("I5").Value = -0.0036992862934230718
$ws.Range("K5").Value = 0.019151276642876192
$ws.Range("B6").Value = 0.32556187564985672
$ws.Range("C6").Value = -0.045902448210449559
$ws.Range("E6").Value = 0.00030289419240961911
$ws.Range("F6").Value = -0.005413571518924771
$ws.Range("H6").Value = -0.010254911839014778
$ws.Range("I6").Value = 0.0002043658735635938
$ws.Range("K6").Value = 0.00030728794769668832
$ws.Range("B7").Value = 0.35656615837284239
$ws.Range("D7").Value = 0.088896553010014437
$ws.Range("E7").Value = -0.021994409515656458
$ws.Range("F7").Value = -0.017088266630824307
$ws.Range("G7").Value = 0.004627591089781221
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = -0.023437185230329227

# Column width adjustments (D:J) reflecting the revised table layout
$ws.Columns("D").ColumnWidth = 13.333333333333332
$ws.Columns("E").ColumnWidth = 14.833333333333332
$ws.Columns("F").ColumnWidth = 14.833333333333332
$ws.Columns("G").ColumnWidth = 14.333333333333332
$ws.Columns("H").ColumnWidth = 14.833333333333332
$ws.Columns("I").ColumnWidth = 14.833333333333332
$ws.Columns("J").ColumnWidth = 14.166666666666666
